$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing data row (row 2): SKU Code and Location change,
#     and its barcode hyperlink needs to point at the new URL ---
$ws.Range("B2").Value = "NC-NTH-8009"
$ws.Range("C2").Value = "R1"

$newUrl2 = "https://barcode.tec-it.com/barcode.ashx?data=OM-NC-NTH-8009-R1"
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), $newUrl2)
$ws.Range("D2").Value = $newUrl2
$ws.Range("D2").Style = "Hyperlink"

# --- Add a new data row (row 3) for the additional SKU ---
$ws.Range("A3").Value = "DIPSHI"
$ws.Range("B3").Value = "NC-NK-9005"
$ws.Range("C3").Value = "R2"

$newUrl3 = "https://barcode.tec-it.com/barcode.ashx?data=DIPSHI-NC-NK-9005-R2"
$ws.Hyperlinks.Add($ws.Range("D3"), $newUrl3)
$ws.Range("D3").Style = "Hyperlink"

Write-Host "Applied SKU/location updates and appended new row with barcode hyperlink."
